# Update column G ("K" strikeout count) values for rows 2-15 on Sheet1.
# This reflects a regeneration of save_data that uses K (strikeouts) instead
# of Strike# for this column, with recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 4
    7  = 0
    8  = 2
    9  = 5
    10 = 5
    11 = 2
    12 = 5
    13 = 4
    14 = 1
    15 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}

$wb.Save()
